# 26 Mayıs verileri eklendi
# Append the 2020-05-26 COVID-19 Turkey data row to the "data" table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet's table (Table3) currently spans A1:E75 (header + 74 data rows).
# Adding a ListRow grows the table (and its autoFilter range) by one row,
# mirroring what Excel does when a user types a new row right under a table.
$lo = $ws.ListObjects.Item(1)
$newRow = $lo.ListRows.Add()

# Fill in the newly added row (row 76) with the 2020-05-26 figures:
# date, test, case, death, recovered
$ws.Cells.Item(76, 1).Value = 43977
$ws.Cells.Item(76, 2).Value = 19853
$ws.Cells.Item(76, 3).Value = 948
$ws.Cells.Item(76, 4).Value = 28
$ws.Cells.Item(76, 5).Value = 1492

# Move the selection to where the user would land after entering the row.
$ws.Range("B77").Select()
